# Add new Spanish/English vocab rows to the "warm" sheet, then move the
# active-tab/selection/view state to match (the edit was an upload that
# appended rows and left the "warm" sheet scrolled/selected near the bottom
# of the new data, instead of the "pictures" sheet that was active before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "warm"

# Row used as the formatting template for column B (style carries the
# plain Calibri font used throughout this block of the table).
$fmtSrc = $ws.Cells.Item(302, 2)

$rows = @(
  @("mas o menos", "pretty much/more or less", "all"),
  @("No va a ser facil venderlo", "it's not going to be easy to sell it", "all"),
  @("se van a enfadar", "they're going to get ungry", "all"),
  @("te vas a poner nervioso", "You're going to get nervous", "all"),
  @("cuanto crees que vamos a tardar?", "How long do you think it's going to take us?", "all"),
  @("Cuanto tardarían los leones en comerme?", "How long would it take the lions to eat me?", "tardé, se tarda/it took me, It takes"),
  @("Cuanto duraría en una habitación llena de leones?", "How long would I last in a room full of lions?", "to last/durar"),
  @("Cuanto crees que vamos a durar?", "How long do you think we're going to last?", "to last/durar"),
  @("Creo que puedo durar 2 minutos bajo el agua", "I think I can last 2 minutes underwater", "to last/durar")
)

$r = 303
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]

  # Match the existing formatting (style index) used in column B.
  $fmtSrc.Copy()
  $ws.Cells.Item($r, 2).PasteSpecial(-4122)

  $r = $r + 1
}

# The "warm" sheet becomes the active tab/sheet (previously "pictures" was
# active), scrolled down near the newly-added rows with the last new cell
# selected and zoomed out a bit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 268
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 73
$ws.Range("C311").Select()
